$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.321.16'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '2.907.47'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '352.54'
$ws.Range("D6").Value = '112.05'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.631'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '39.93'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  +3.21%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '19.86'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.80'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '3.364.23'
$ws.Range("E15").Value = '  +3.63%  '
$ws.Range("E16").Value = '  +6.02%  '
$ws.Range("D17").Value = '2.904.95'
$ws.Range("E17").Value = '  +3.39%  '
$ws.Range("D18").Value = '52.342.20'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = '7.62'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("D21").Value = '14.19'
$ws.Range("E21").Value = '  +4.20%  '
$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").Value = '70.95'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = '269.93'
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '0.173'
$ws.Range("E26").Value = '  +6.97%  '
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  +2.67%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '6.65'
$ws.Range("E30").Value = '  +8.29%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '37.93'
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = '6.36'
$ws.Range("E32").Value = '  +12.51%  '
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").Value = '0.0978'
$ws.Range("E34").Value = '  +10.52%  '
$ws.Range("D35").Value = '53.32'
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").Value = '0.0451'
$ws.Range("E36").Value = '  +1.73%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("D39").Value = '18.82'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("E40").Value = '  +3.18%  '
$ws.Range("D41").Value = '2.84'
$ws.Range("E41").Value = '  +13.87%  '
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = '23.34'
$ws.Range("E43").Value = '  +6.23%  '
$ws.Range("D44").Value = '121.38'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").Value = '2.61'
$ws.Range("E45").Value = '  +7.82%  '
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").Value = '3.55'
$ws.Range("E47").Value = '  +3.98%  '
$ws.Range("D48").Value = '2.196.78'
$ws.Range("E48").Value = '  +4.00%  '
$ws.Range("D49").Value = '0.264'
$ws.Range("E49").Value = '  +21.60%  '
$ws.Range("E50").Value = '  +9.50%  '
$ws.Range("D51").Value = '0.971'
$ws.Range("E51").Value = '  +2.55%  '
